$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: A:D become uniform (closest achievable to 19.85546875 via the
# character-width ColumnWidth property, which snaps to a 1/6-pixel grid) ---
$ws.Columns("A:D").ColumnWidth = 19

# --- Row 1 (header) height: 45 -> 30 ---
$ws.Rows(1).RowHeight = 30

# --- Header row A1:C1 (font18, fill, center/top/wrap, no border): vertical top -> center ---
$ws.Range("A1:C1").VerticalAlignment = -4108

# --- Data row A2:C2 (font18, fill, border, center/top/wrap): vertical top -> center ---
$ws.Range("A2:C2").VerticalAlignment = -4108

# --- D1 / D2 (fill, vertical-top + wrap, no horizontal) become horizontal=center,
# vertical=center, wrap. Build the target format once on a scratch cell (seeded from
# A3's "wrap only" style, which already matches on every attribute except alignment)
# and then copy that exact format onto D1 and D2, avoiding throw-away style entries. ---
$scratch = $ws.Range("Z1")
$ws.Range("A3").Copy()
$scratch.PasteSpecial(-4122)
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$scratch.Clear()

# --- New data values for row 2 (grand total calc) ---
$ws.Range("A2").Value = 98
$ws.Range("B2").Value = 254
$ws.Range("C2").Value = 17
$ws.Range("D2").Value = 369

# --- Selection moves from B2 to A2 ---
$ws.Range("A2").Select()

Write-Host "done"
